$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tracked-stock table (rows 11-20) is rewritten in place: "신세계", "KB금융" and
# "메리츠화재" drop out of the list, the remaining names ripple up, and three newly
# tracked stocks (휠라홀딩스, 현대건설기계, 대한항공) are appended at the bottom.

$ws.Range("A11").Value = 44117
$ws.Range("B11").Value = "SK텔레콤"
$ws.Range("C11").Value = 240000

$ws.Range("A12").Value = 44117
$ws.Range("B12").Value = "테스"
$ws.Range("C12").Value = 24700

$ws.Range("A13").Value = 44119
$ws.Range("B13").Value = "케이씨텍"
$ws.Range("C13").Value = 23950

$ws.Range("A14").Value = 44120
$ws.Range("B14").Value = "한온시스템"
$ws.Range("C14").Value = 12550

$ws.Range("A15").Value = 44123
$ws.Range("B15").Value = "풍산"
$ws.Range("C15").Value = 25950

$ws.Range("A16").Value = 44124
$ws.Range("B16").Value = "삼성전자"
$ws.Range("C16").Value = 59500

$ws.Range("A17").Value = 44124
$ws.Range("B17").Value = "NAVER"
$ws.Range("C17").Value = 308875

$ws.Range("A18").Value = 44125
$ws.Range("B18").Value = "휠라홀딩스"
$ws.Range("C18").Value = 40250

$ws.Range("A19").Value = 44126
$ws.Range("B19").Value = "현대건설기계"
$ws.Range("C19").Value = 30000

$ws.Range("A20").Value = 44127
$ws.Range("B20").Value = "대한항공"
$ws.Range("C20").Value = 21300

# Row 11 (SK텔레콤) and row 17 (NAVER) are the "touched target price" rows, so they
# keep the highlighted name style; every other name cell uses the normal style.
# Pull matching style examples from rows that already carry the right look.
$ws.Range("B17").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# The last row's recommendation date uses a plain mm-dd-yy short-date format instead
# of the Korean long-date format used for the other rows in column A.
$ws.Range("A20").NumberFormat = "mm-dd-yy"

$ws.CutCopyMode = 0

# Match the saved workbook's active selection.
$ws.Range("E11").Select()
